# Refresh the cryptocurrency price/volume table (GitHub Actions data update).
# Column D ("Price") cells are stored as plain text in the workbook even when
# their content looks numeric, so force a text number-format while writing them
# and then restore the cell to its original, style-less "Normal" appearance.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '67.849.80'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E2').Value = '  -0.34%  '

$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '3.802.34'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E3').Value = '  -1.84%  '

$ws.Range('E4').Value = '  +0.10%  '

$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '599.34'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  -0.08%  '

$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '168.59'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  +0.79%  '

$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '3.797.01'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E7').Value = '  -1.89%  '

$ws.Range('E8').Value = '  +0.04%  '

$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.529'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E9').Value = '  +0.40%  '

$ws.Range('E10').Value = '  -0.08%  '

$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '6.51'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  +1.10%  '

$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '0.462'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E12').Value = '  +1.00%  '

$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '0.0000276'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E13').Value = '  +10.97%  '

$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '36.89'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  -0.19%  '

$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '4.444.23'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  -1.63%  '

$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '3.784.24'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E16').Value = '  -2.25%  '

$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '68.006.31'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E17').Value = '  -0.16%  '

$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '18.36'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E18').Value = '  +0.70%  '

$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '7.47'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  +1.02%  '

$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '0.112'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  +0.42%  '

$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '10.85'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  -0.19%  '

$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '469.79'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  +0.81%  '

$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '0.735'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E23').Value = '  +0.41%  '

$ws.Range('E24').Value = '  -7.72%  '

$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '83.43'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  -0.05%  '

$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '2.31'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  +3.13%  '

$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '12.21'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  +0.80%  '

$ws.Range('E28').Value = '  +3.30%  '

$ws.Range('E29').Value = '  -0.01%  '

$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '2.92'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  -1.04%  '

$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '3.953.52'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  -1.71%  '

$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '7.73'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E32').Value = '  -1.96%  '

$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '2.28'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  -1.08%  '

$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '30.81'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E34').Value = '  -1.45%  '

$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '9.32'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  -0.41%  '

$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '3.771.89'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E36').Value = '  -1.98%  '

$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '0.107'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E37').Value = '  +2.60%  '

$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '3.79'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E38').Value = '  +9.78%  '

$ws.Range('B39').Value = 'Filecoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '5.97'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  +1.10%  '

$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '0.140'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  -0.08%  '

$ws.Range('B41').Value = 'Mantle'
$ws.Range('C41').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '1.01'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  -1.30%  '

$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  +0.00%  '

$ws.Range('E43').Value = '  +1.65%  '

$ws.Range('E44').Value = '  +0.01%  '

$ws.Range('B45').Value = 'Cosmos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '8.77'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  +2.58%  '

$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '1.97'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E46').Value = '  -0.45%  '

$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '46.44'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  -1.71%  '

$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '407.05'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  -5.94%  '

$ws.Range('B49').Value = 'FLOKI'
$ws.Range('C49').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '0.000286'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  +2.56%  '

$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '0.0360'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  +0.80%  '

$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '141.91'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  -1.54%  '
